# pmSheet.xlsx - Battle System Core Mechanics update
# - Implement the Battle System: Effort (D3) 9 -> 12 (Remaining E3 recalculates)
# - > Implement Core Mechanics (Fight, Run): Effort (D13) 2 -> 3 (Remaining E13 recalculates)
# - > Implement Shield Mechanic: Effort (D14) 14 -> 15 (Remaining E14 recalculates)
# - Active selection moved to D15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 12
$ws.Range("D13").Value = 3
$ws.Range("D14").Value = 15

$ws.Range("D15").Select()
